$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers (losing significant trailing zeros), matching
# the source data which stores these as plain text strings.
$textCells = @("D5", "D6", "D8", "D14", "D19", "D20", "D21", "D22", "D24", "D27", "D29", "D31", "D35", "D39", "D40", "D43", "D44", "D46", "D47", "D48", "D51", "D41", "D42")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = '60.429.58'
$ws.Range("E2").Value = '  +3.69%  '
$ws.Range("D3").Value = '2.323.79'
$ws.Range("E3").Value = '  +1.52%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '545.28'
$ws.Range("E5").Value = '  +1.82%  '
$ws.Range("D6").Value = '131.03'
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '0.580'
$ws.Range("E8").Value = '  -1.34%  '
$ws.Range("D9").Value = '2.321.84'
$ws.Range("E9").Value = '  +1.68%  '
$ws.Range("E10").Value = '  +0.69%  '
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("E13").Value = '  +0.40%  '
$ws.Range("D14").Value = '23.68'
$ws.Range("E14").Value = '  +0.23%  '
$ws.Range("D15").Value = '60.426.45'
$ws.Range("E15").Value = '  +3.81%  '
$ws.Range("D16").Value = '2.740.92'
$ws.Range("E16").Value = '  +1.66%  '
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("D18").Value = '2.321.23'
$ws.Range("E18").Value = '  +0.98%  '
$ws.Range("D19").Value = '10.60'
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("D20").Value = '4.14'
$ws.Range("E20").Value = '  -1.08%  '
$ws.Range("D21").Value = '315.31'
$ws.Range("E21").Value = '  +0.54%  '
$ws.Range("D22").Value = '6.62'
$ws.Range("E22").Value = '  +1.24%  '
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").Value = '64.00'
$ws.Range("E24").Value = '  +1.20%  '
$ws.Range("E25").Value = '  +1.51%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").Value = '7.85'
$ws.Range("E27").Value = '  -1.68%  '
$ws.Range("E28").Value = '  +4.32%  '
$ws.Range("D29").Value = '173.72'
$ws.Range("E29").Value = '  +1.80%  '
$ws.Range("E30").Value = '  +9.62%  '
$ws.Range("D31").Value = '1.73'
$ws.Range("E31").Value = '  +1.66%  '
$ws.Range("D32").Value = '0.0₃0731'
$ws.Range("E32").Value = '  +0.73%  '
$ws.Range("E33").Value = '  +2.04%  '
$ws.Range("E34").Value = '  +11.09%  '
$ws.Range("D35").Value = '0.380'
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("D39").Value = '4.05'
$ws.Range("E39").Value = '  +3.41%  '
$ws.Range("D40").Value = '325.89'
$ws.Range("E40").Value = '  +12.66%  '
$ws.Range("D43").Value = '138.30'
$ws.Range("E43").Value = '  -1.74%  '
$ws.Range("D44").Value = '3.48'
$ws.Range("E44").Value = '  +1.16%  '
$ws.Range("E45").Value = '  -0.92%  '
$ws.Range("D46").Value = '19.25'
$ws.Range("E46").Value = '  +6.24%  '
$ws.Range("D47").Value = '0.0496'
$ws.Range("E47").Value = '  +0.43%  '
$ws.Range("D48").Value = '0.560'
$ws.Range("E48").Value = '  +1.17%  '
$ws.Range("E49").Value = '  +1.03%  '
$ws.Range("D50").Value = '0.0₆0215'
$ws.Range("E50").Value = '  +19.01%  '
$ws.Range("D51").Value = '11.02'
$ws.Range("E51").Value = '  +0.60%  '
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").Value = '37.93'
$ws.Range("E41").Value = '  -1.03%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '1.53'
$ws.Range("E42").Value = '  +2.11%  '
